$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A92").Value = "G1"
$ws.Range("B92").Value = "Test1"
$ws.Range("C92").Value = 45906
$ws.Range("C92").NumberFormat = "YYYY-MM-DD"
$ws.Range("D92").Value = 0.6454454647550686
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = -0.01

$ws.Range("A93").Value = "G2"
$ws.Range("B93").Value = "sedrftgyhuioygtfrd"
$ws.Range("C93").Value = 45906
$ws.Range("C93").NumberFormat = "YYYY-MM-DD"
$ws.Range("D93").Value = 0.6454454647550686
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = -0.01
